$d = $word.ActiveDocument

# --- Update the four "Group Members" paragraphs ---
$d.Content.Find.Execute("Yousif Manhal Talal 20193866", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Yousif ", 2)

$d.Content.Find.Execute("Mahdi AbdulHussain Hasan 202009297", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Adnan", 2)

$d.Content.Find.Execute("Mahmood Almajed 202006385", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Hassan", 2)

$d.Content.Find.Execute("Abdullah Ameen Naji 20175012", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Salman", 2)

# --- Update the "Name" column of the use-case table ---
$table = $d.Tables.Item(1)

# Row 2 ("Yousif" / "Manhal") -> "Yousif "
$cell = $table.Cell(2, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$rng.Text = "Yousif "

# Row 3 ("Abdulla" / " Ameen") -> "Adnan"
$cell = $table.Cell(3, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$rng.Text = "Adnan"

# Row 4 ("Mahdi" / " Abdul" / "Hussain") -> "Hassan"
$cell = $table.Cell(4, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$rng.Text = "Hassan"

# Row 5 ("Mahmood" / " Almajed ") -> "Salman" / " "
$cell = $table.Cell(5, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$rng.Text = "Salman "
